$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record a homework grade of 5 for the student in row 31 (Федченко Кирилл)
# across the four assignment columns C:F, matching the pattern already
# used for the other graded students (e.g. row 6).
$ws.Range("C31").Value = 5
$ws.Range("D31").Value = 5
$ws.Range("E31").Value = 5
$ws.Range("F31").Value = 5

# Reflect that the user scrolled down to / ended up with F31 selected in
# the bottom-right (unfrozen) pane after entering the grades.
$av = $excel.ActiveWindow
$av.ScrollRow = 13
$av.ScrollColumn = 3
$ws.Range("F31").Select()
